$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.859.23"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.814.32"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'309.39"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4325"
$ws.Range("E7").Value = "  +2.08%  "
$ws.Range("D8").Value = "'0.3711"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("D9").Value = "'0.07281"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'0.8683"
$ws.Range("E10").Value = "  +2.67%  "
$ws.Range("D11").Value = "'20.97"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").Value = "1.959.66"
$ws.Range("E12").Value = "  +9.61%  "
$ws.Range("D13").Value = "'6.649"
$ws.Range("E13").Value = "  +4.17%  "
$ws.Range("D14").Value = "'5.365"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("D15").Value = "'0.06925"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "'80.72"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "'0.000008927"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "'15.30"
$ws.Range("E20").Value = "  +1.52%  "
$ws.Range("D21").Value = "26.883.35"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'5.216"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "'11.20"
$ws.Range("E23").Value = "  +0.77%  "
$ws.Range("D24").Value = "2.183.82"
$ws.Range("E24").Value = "  +8.26%  "
$ws.Range("D25").Value = "'154.00"
$ws.Range("E25").Value = "  +0.56%  "
$ws.Range("D26").Value = "'1.871"
$ws.Range("E26").Value = "  -4.32%  "
$ws.Range("D27").Value = "'18.28"
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "'5.245"
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "'1.896"
$ws.Range("E29").Value = "  +14.16%  "
$ws.Range("D30").Value = "'115.30"
$ws.Range("E30").Value = "  +1.01%  "
$ws.Range("D31").Value = "'0.08953"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "'0.7590"
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("D33").Value = "'1.175"
$ws.Range("E33").Value = "  +6.87%  "
$ws.Range("D34").Value = "'4.449"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "'2.806"
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("D36").Value = "'1.007"
$ws.Range("D37").Value = "'1.131"
$ws.Range("E37").Value = "  +4.59%  "
$ws.Range("D38").Value = "'0.05244"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "'0.5093"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").Value = "'0.1653"
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("D42").Value = "'2.665"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "'6.568"
$ws.Range("E43").Value = "  +10.16%  "
$ws.Range("D44").Value = "'8.295"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'106.98"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").Value = "'10.45"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'1.657"
$ws.Range("E48").Value = "  +3.17%  "
$ws.Range("D49").Value = "'0.4586"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Value = "'0.06294"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("D51").Value = "'1.819"
$ws.Range("E51").Value = "  +5.30%  "
